$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.049741383723816215
$ws.Range("B1").Value = -0.049741383734224313

$ws.Range("A2").Value = 0.022146649361710815
$ws.Range("B2").Value = -0.022146649378215099

$ws.Range("A3").Value = -0.03578125063916851
$ws.Range("B3").Value = 0.035781250611754446

$ws.Range("A4").Value = -0.017750536747277772
$ws.Range("B4").Value = 0.017750536705936636
